$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 54 values (B54 and D54 recalculated)
$ws.Cells.Item(54, 2).Value = 160
$ws.Cells.Item(54, 4).Value = 74.07407407407408

# Add new row 55 for 2025-03
$ws.Cells.Item(55, 1).Value = "2025-03"
$ws.Cells.Item(55, 2).Value = 77
$ws.Cells.Item(55, 3).Value = 231
$ws.Cells.Item(55, 4).Value = 33.33333333333333
